$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header for column E changes from "Total NH4 Mass (g)" to "Total NH3 Mass (g)"
# (the protocol now tracks NH3 instead of NH4).
$ws.Range("E1").Value = "Total NH3 Mass (g)"

# The mass values feeding the yield calc are now scaled by a 0.75 recovery
# factor on the first (measured) term before adding the second constant.
$ws.Range("E2").Formula = "=0.75*0.127+0.441"
$ws.Range("E3").Formula = "=0.75*0.067308+0.234"
$ws.Range("E4").Formula = "=0.75*0.157+0.546"
$ws.Range("E5").Formula = "=0.75*0.09293+0.323"

# Leave the selection on the updated header cell, matching the saved view.
[void]$ws.Range("E1").Select()
